# 2021 Q2 Quarterly Report WIP
#
# Two kinds of edits to the stats table:
#   1. Every run that used "Helvetica" for hAnsi/eastAsia/cs (while ascii was
#      already "Times New Roman") gets normalized to "Times New Roman" across
#      the board.
#   2. A handful of count/percentage table cells get updated numbers.

$d = $word.ActiveDocument

# --- 1. Font normalization: Helvetica -> Times New Roman (hAnsi/eastAsia/cs) ---
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Font.NameFarEast = "Helvetica"
$find.Replacement.Font.NameAscii = "Times New Roman"
$find.Replacement.Font.Name = "Times New Roman"
$find.Replacement.Font.NameFarEast = "Times New Roman"
$find.Replacement.Font.NameBi = "Times New Roman"
$find.Execute("", $false, $false, $false, $false, $false, $true, 1, $true, "", 2) | Out-Null

# --- 2. Update the statistic cell values in-place (row, col are 1-based) ---
#     Using Cell(row,col).Range keeps the run's xml:space="preserve" marker
#     (a plain Find/Replace on text drops it once there's no leading/trailing
#     whitespace left to protect).
$tbl = $d.Tables.Item(1)

$cellReplacements = @(
    @(2, 2, "196 (71.5)"),
    @(2, 3, "69 (74.2)"),
    @(2, 4, "61 (66.3)"),
    @(2, 5, "66 (74.2)"),
    @(3, 2, "32 (11.7)"),
    @(3, 3, "13 (14.0)"),
    @(3, 4, "11 (12.0)"),
    @(3, 5, "8 (9.0)"),
    @(4, 2, "130 (47.4)"),
    @(4, 3, "43 (46.2)"),
    @(4, 4, "40 (43.5)"),
    @(4, 5, "47 (52.8)"),
    @(5, 2, "62 (22.6)"),
    @(5, 3, "18 (19.4)"),
    @(5, 4, "26 (28.3)"),
    @(5, 5, "18 (20.2)")
)

foreach ($entry in $cellReplacements) {
    $row = $entry[0]
    $col = $entry[1]
    $newText = $entry[2]
    $cellRange = $tbl.Cell($row, $col).Range
    $cellRange.MoveEnd(1, -1) | Out-Null
    $cellRange.Text = $newText
}
